# Applies the capital-structure database refresh for the France "Bank (Money Center)"
# rows (2-5): updated growth/margin/debt metrics, and a handful of now-obsolete
# buyback/growth cells that are cleared or newly populated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = 0.00396
$ws.Range("E2").Value = 0.02349
$ws.Range("F2").Value = 0.0233
$ws.Range("K2").Value = 13479.3
$ws.Range("L2").Value = 0.1472342981977062
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 564901.4
$ws.Range("V2").Value = 4.710070805137458
$ws.Range("W2").Value = 0.07310158026835856
$ws.Range("X2").Value = 0.2703759066830965
$ws.Range("Y2").Value = -0.197274326414738
$ws.Range("Z2").Value = 0.06874464469456164
$ws.Range("AB2").Value = 0.03540300590080994
$ws.Range("AC2").Value = -0.03540300590080994
$ws.Range("AD2").Value = 1489519.7
$ws.Range("AF2").Value = 1489519.7
$ws.Range("AG2").Value = 924618.3000000002
$ws.Range("AH2").Value = 0.9254810869148522
$ws.Range("AI2").Value = 0.8331166158786404
$ws.Range("AJ2").Value = 0.8851807533767312
$ws.Range("AK2").Value = 0.7560326878318923

# Row 3
$ws.Range("D3").Value = 0.00396
$ws.Range("E3").Value = -0.00222
$ws.Range("F3").Value = 0.0067
$ws.Range("K3").Value = 8600.299999999999
$ws.Range("L3").Value = 0.1816382778757511
$ws.Range("U3").Value = 370896
$ws.Range("V3").Value = 5.631488086290693
$ws.Range("W3").Value = 0.0736054185391888
$ws.Range("X3").Value = 0.2084047404425275
$ws.Range("Y3").Value = -0.1347993219033387
$ws.Range("Z3").Value = 0.06970314059767588
$ws.Range("AB3").Value = 0.03537416153040585
$ws.Range("AC3").Value = -0.03537416153040585
$ws.Range("AD3").Value = 639193.3
$ws.Range("AF3").Value = 639193.3
$ws.Range("AG3").Value = 268297.3
$ws.Range("AH3").Value = 0.9065872080225299
$ws.Range("AI3").Value = 0.823772026176945
$ws.Range("AJ3").Value = 0.8029045506562158
$ws.Range("AK3").Value = 0.6623991732148063

# Row 4
$ws.Range("D4").Value = 0.04219999999999999
$ws.Range("E4").Value = 0.0492
$ws.Range("F4").Value = 0.0233
$ws.Range("K4").Value = 4965.9
$ws.Range("L4").Value = 0.2347854454677837
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.07310158026835856
$ws.Range("X4").Value = 0.2703759066830965
$ws.Range("Y4").Value = -0.197274326414738
$ws.Range("Z4").Value = 0.05476397490118259
$ws.Range("AB4").Value = 0.03540300590080994
$ws.Range("AC4").Value = -0.03540300590080994
$ws.Range("AD4").Value = 478721.2
$ws.Range("AF4").Value = 478721.2
$ws.Range("AG4").Value = 478721.2
$ws.Range("AH4").Value = 0.9293694933364809
$ws.Range("AI4").Value = 0.8500200909784849
$ws.Range("AJ4").Value = 0.9293694933364809
$ws.Range("AK4").Value = 0.8500200909784849

# Row 5
$ws.Range("D5").Value = -0.0363
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 0.321
$ws.Range("K5").Value = -86.90000000000001
$ws.Range("L5").Value = -0.003769950587183903
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 194005.4
$ws.Range("V5").Value = 10.96589926349644
$ws.Range("W5").Value = -0.001250820446092374
$ws.Range("X5").Value = 0.4111932905859613
$ws.Range("Y5").Value = -0.1347993219033387
$ws.Range("Z5").Value = 0.06970314059767588
$ws.Range("AB5").Value = 0.03619853991575021
$ws.Range("AC5").Value = -0.03619853991575021
$ws.Range("AD5").Value = 371605.2
$ws.Range("AF5").Value = 371605.2
$ws.Range("AG5").Value = 177599.8
$ws.Range("AH5").Value = 0.9545547370143456
$ws.Range("AI5").Value = 0.8280604270646292
$ws.Range("AJ5").Value = 0.8029045506562158
$ws.Range("AK5").Value = 0.6623991732148063
